$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "nama_kategori"
$ws.Range("E5").Select()
